$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.503.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.11%  '

$ws.Range("D3").Value = "'3.436.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.55%  '

$ws.Range("D5").Value = "'592.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.03%  '

$ws.Range("D6").Value = "'135.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.87%  '

$ws.Range("D7").Value = "'3.435.80"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.57%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").Value = "'0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.34%  '

$ws.Range("E10").Value = '  -6.43%  '

$ws.Range("E11").Value = '  -9.25%  '

$ws.Range("D12").Value = "'0.378"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.00%  '

$ws.Range("D13").Value = "'4.016.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.69%  '

$ws.Range("D14").Value = "'0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -11.20%  '

$ws.Range("D15").Value = "'26.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -9.58%  '

$ws.Range("D16").Value = "'3.452.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.17%  '

$ws.Range("D17").Value = "'65.508.15"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = '  -2.31%  '

$ws.Range("D19").Value = "'9.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.77%  '

$ws.Range("D20").Value = "'5.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.38%  '

$ws.Range("D21").Value = "'13.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.63%  '

$ws.Range("D22").Value = "'394.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.17%  '

$ws.Range("D23").Value = "'0.546"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.20%  '

$ws.Range("D24").Value = "'73.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.72%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").Value = "'3.578.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.38%  '

$ws.Range("D27").Value = "'0.0000105"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.80%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = "'7.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.58%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.93%  '

$ws.Range("D31").Value = "'8.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.81%  '

$ws.Range("D32").Value = "'3.441.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.31%  '

$ws.Range("E33").Value = '  -0.05%  '

$ws.Range("D34").Value = "'0.147"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.40%  '

$ws.Range("D35").Value = "'22.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.28%  '

$ws.Range("D36").Value = "'171.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.12%  '

$ws.Range("D37").Value = "'1.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -13.19%  '

$ws.Range("D38").Value = "'6.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.15%  '

$ws.Range("E39").Value = '  -7.29%  '

$ws.Range("D40").Value = "'4.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -11.55%  '

$ws.Range("D41").Value = "'0.0771"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.44%  '

$ws.Range("D42").Value = "'0.822"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.96%  '

$ws.Range("D43").Value = "'43.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.99%  '

$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("D45").Value = "'4.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -14.40%  '

$ws.Range("D46").Value = "'1.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -12.28%  '

$ws.Range("D47").Value = "'1.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.26%  '

$ws.Range("D48").Value = "'22.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.43%  '

$ws.Range("D49").Value = "'6.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.86%  '

$ws.Range("D50").Value = "'2.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -15.86%  '

$ws.Range("D51").Value = "'2.200.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.55%  '
